$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 2, pushing existing rows 2-4 down to 3-5.
$ws.Rows("2:2").Insert()

# Populate the new row 2 with the slug-style identifiers that link the
# header row (row 1) to the existing metadata rows (now rows 3-5), so two
# columns can be related to build hierarchical SKOS concepts.
$ws.Range("A2").Value = "edad"
$ws.Range("B2").Value = "estado-civil"
$ws.Range("C2").Value = "personas-residentes-viviendas-familiares"
$ws.Range("D2").Value = "comarca-nombre"
$ws.Range("E2").Value = "comarca-codigo"
$ws.Range("F2").Value = "provincia-codigo"
$ws.Range("G2").Value = "aragon"
$ws.Range("H2").Value = "provincia-nombre"
$ws.Range("I2").Value = "sexo"
